# ADD results from server
# Update row 2 values (B2:Q2) with new results from server

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 37.84278171914715
$ws.Range("C2").Value = 38.74047329475788
$ws.Range("D2").Value = 36.93294442240584
$ws.Range("E2").Value = 37.84294501151046
$ws.Range("F2").Value = 37.94553809617636
$ws.Range("G2").Value = 37.18388087321011
$ws.Range("H2").Value = 40.66793297706615
$ws.Range("I2").Value = 33.84811172354475
$ws.Range("J2").Value = 38.54953366014051
$ws.Range("K2").Value = 37.90966494870825
$ws.Range("L2").Value = 37.84710463969298
$ws.Range("M2").Value = 37.05173872635862
$ws.Range("N2").Value = 19.84180137267095
$ws.Range("O2").Value = 33.16907792640307
$ws.Range("P2").Value = 41.39802372017177
$ws.Range("Q2").Value = 33.213818479076
